$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Date Colours")

# ---------------------------------------------------------------------
# Sheet1 / Table1: add a new "Isolated" column and a new data row
# ---------------------------------------------------------------------
$lo = $ws1.ListObjects.Item(1)

# Add the 8th table column (this becomes column H) and name it "Isolated"
$newCol = $lo.ListColumns.Add()
$ws1.Range("H1").Value = "Isolated"

# Give column H roughly the same visual width Excel produced (bestFit-ish)
$ws1.Columns.Item(8).ColumnWidth = 9.5

# Fill in the "Isolated" / "Not isolated" values for the existing rows
$ws1.Range("H2").Value = "Not isolated"
$ws1.Range("H3").Value = "Not isolated"
$ws1.Range("H4").Value = "Not isolated"
$ws1.Range("H5").Value = "Not isolated"
$ws1.Range("H6").Value = "Not isolated"
$ws1.Range("H7").Value = "Not isolated"
$ws1.Range("H8").Value = "Not isolated"
$ws1.Range("H9").Value = "Not isolated"
$ws1.Range("H10").Value = "Not isolated"
$ws1.Range("H11").Value = "Not isolated"
$ws1.Range("H12").Value = "Not isolated"
$ws1.Range("H13").Value = "Isolated"
$ws1.Range("H14").Value = "Isolated"
$ws1.Range("H15").Value = "Isolated"
$ws1.Range("H16").Value = "Isolated"
$ws1.Range("H17").Value = "Isolated"

# ---------------------------------------------------------------------
# "Date Colours" sheet: shift the colour-gradient values down one row,
# introducing a new lighter shade at the top of the gradient
# ---------------------------------------------------------------------
$ws2.Range("B2").Value = "#ffeee7"
$ws2.Range("B3").Value = "#ffddcf"
$ws2.Range("B4").Value = "#fdccb8"
$ws2.Range("B5").Value = "#fabba1"
$ws2.Range("B6").Value = "#f6ab8b"
$ws2.Range("B7").Value = "#f19a75"
$ws2.Range("B8").Value = "#eb895f"

# Add a new table row (row 18) and populate it
$newRow = $lo.ListRows.Add()
$ws1.Range("A18").Value = 44379
$ws1.Range("A18").NumberFormat = $ws1.Range("A17").NumberFormat()
$ws1.Range("B18").Value = "S1 m"
$ws1.Range("C18").Value = "S6 child"
$ws1.Range("D18").Value = "South Australia"
$ws1.Range("F18").Value = "Household"
$ws1.Range("G18").Value = "Delta (B.1.617.2)"
$ws1.Range("H18").Value = "Isolated"

# ---------------------------------------------------------------------
# View state: Sheet1 becomes the active / selected tab, with H17
# selected; the "Date Colours" sheet selection moves to F2:L2
# ---------------------------------------------------------------------
$ws2.Range("F2:L2").Select()
$ws1.Activate()
$ws1.Range("H17").Select()
